$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = -7.614799999999992
$ws.Range("C12").Value = -11.3907
$ws.Range("D23").Value = -7.860600000000002
$ws.Range("C27").Value = -12.38239999999999
$ws.Range("D28").Value = -7.974099999999996
$ws.Range("C32").Value = -13.1219
$ws.Range("D32").Value = -8.328899999999997
$ws.Range("D34").Value = -7.802700000000001
$ws.Range("C36").Value = -11.9093
$ws.Range("C38").Value = -12.2841
$ws.Range("D42").Value = -8.693599999999995
$ws.Range("C46").Value = -14.66469999999999
$ws.Range("D49").Value = -8.003300000000001
$ws.Range("C54").Value = -12.4267
$ws.Range("D54").Value = -7.970100000000007
$ws.Range("C55").Value = -13.6874
$ws.Range("C56").Value = -12.90869999999999
$ws.Range("C67").Value = -12.48909999999999
$ws.Range("C69").Value = -11.77859999999999
$ws.Range("C72").Value = -11.75570000000001
$ws.Range("D78").Value = -7.966500000000001
$ws.Range("D80").Value = -7.734100000000001
$ws.Range("C83").Value = -13.4827
$ws.Range("C86").Value = -14.36329999999999
$ws.Range("C91").Value = -12.17599999999999
$ws.Range("C93").Value = -10.5382
$ws.Range("D97").Value = -8.200299999999993
$ws.Range("C99").Value = -12.7435
$ws.Range("D99").Value = -8.124499999999998
$ws.Range("D101").Value = -7.6287
$ws.Range("C104").Value = -12.44310000000001
